# edit.ps1 - PowerShell-style PowerPoint COM-interop script
#
# Applies the two logical changes captured by the target diff:
#
#   1. Three tables (on slides 14, 15 and 16) get their gallery table
#      style switched from the plain "no style, no grid" built-in style
#      ({0BB2755F-F846-47A5-A234-DAD7DB79540F}) to the built-in themed
#      table style {AAA93E62-1964-4839-8EC1-2489DB551044} - exactly what
#      clicking a different style in the Table Design gallery would do.
#
#   2. The deck's theme colour scheme (ppt/theme/theme1.xml, used by the
#      single slide master / all slides) is switched from the custom
#      "Integral" / "Red Violet" palette to the standard Office palette
#      - exactly what picking "Office" from Design > Variants > Colors
#      would do. The font scheme and format scheme are already identical
#      between the two themes in this deck, so only the 12 colour scheme
#      slots need to change.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------
$newTableStyle = "{AAA93E62-1964-4839-8EC1-2489DB551044}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Swap the theme colour scheme to the standard Office palette --
# VBA/COM RGB values are packed as 0x00BBGGRR (red in the low byte), so
# build each value from its hex triplet with a small helper instead of
# relying on a VBA RGB() intrinsic.
function HexToVbRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office theme colours, in MsoThemeColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToVbRgb $officeColors[$i - 1]
}
